# Got the program to shuffle the deck properly.
#
# The blank paragraph immediately before "Card" becomes a new "Game"
# heading, and the _GoBack bookmark (previously sitting at the end of the
# "Effects" paragraph) moves to sit right after the new "Game" text.

$d = $word.ActiveDocument
$CR = [char]13

# Locate the target blank paragraph: the empty paragraph immediately
# preceding the paragraph whose text is "Card".
$targetIndex = -1
for ($i = 1; $i -lt $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $text = $p.Range.Text.TrimEnd($CR)
    $nextText = $d.Paragraphs.Item($i + 1).Range.Text.TrimEnd($CR)
    if ($text -eq "" -and $nextText -eq "Card") {
        $targetIndex = $i
    }
}

$target = $d.Paragraphs.Item($targetIndex)

# Fill it in with "Game" plus one throw-away trailing character. The
# trailing character keeps the bookmark's insertion point from landing
# exactly on the paragraph mark (a collapsed position bookmarks can't
# anchor to in this runtime), and the character gets deleted right after.
$target.Range.Text = "GameX"
$target = $d.Paragraphs.Item($targetIndex)

$bookmarkPos = $target.Range.Start + 4
$bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)

# Adding a bookmark under a name that already exists elsewhere moves it,
# so this both relocates _GoBack off of "Effects" and plants it here.
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

$trailing = $d.Range($bookmarkPos, $bookmarkPos + 1)
$trailing.Delete()
